$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Budgets")
$ws.Activate()

# Fill in the two rows of budget data that were previously left blank (0)
# in column H so the "bad data" test fixture exercises a non-trivial value.
$ws.Cells.Item(2, 8).Value = 1   # H2
$ws.Cells.Item(3, 8).Value = 2   # H3

# Leave the selection where the user last clicked before saving.
$ws.Range("H4").Select()
